$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 186, pushing existing rows 186-199 down to 189-202
$ws.Rows.Item(186).Insert()
$ws.Rows.Item(186).Insert()
$ws.Rows.Item(186).Insert()

# Fill in the 3 new rows (186, 187, 188) with weekly Region del Maule data
# Row 186
$ws.Range("A186").Value = 4
$ws.Range("B186").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C186").Value = "Los Lagos"
$ws.Range("D186").Value = 44578
$ws.Range("E186").Value = 10
$ws.Range("F186").Value = 100112028
$ws.Range("G186").Value = "Sandia"
$ws.Range("H186").Value = "Sin especificar"
$ws.Range("I186").Value = "Primera"
$ws.Range("J186").Value = 2000
$ws.Range("K186").Value = 3000
$ws.Range("L186").Value = 3000
$ws.Range("M186").Value = 3000
$ws.Range("N186").Value = "$/unidad"
$ws.Range("O186").Value = "Región del Maule"
$ws.Range("P186").Value = 3000
$ws.Range("Q186").Value = 1
$ws.Range("R186").Value = "Hortaliza"

# Row 187
$ws.Range("A187").Value = 4
$ws.Range("B187").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C187").Value = "Los Lagos"
$ws.Range("D187").Value = 44578
$ws.Range("E187").Value = 10
$ws.Range("F187").Value = 100112028
$ws.Range("G187").Value = "Sandia"
$ws.Range("H187").Value = "Sin especificar"
$ws.Range("I187").Value = "Segunda"
$ws.Range("J187").Value = 2000
$ws.Range("K187").Value = 2500
$ws.Range("L187").Value = 2500
$ws.Range("M187").Value = 2500
$ws.Range("N187").Value = "$/unidad"
$ws.Range("O187").Value = "Región del Maule"
$ws.Range("P187").Value = 2500
$ws.Range("Q187").Value = 1
$ws.Range("R187").Value = "Hortaliza"

# Row 188
$ws.Range("A188").Value = 4
$ws.Range("B188").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C188").Value = "Los Lagos"
$ws.Range("D188").Value = 44578
$ws.Range("E188").Value = 10
$ws.Range("F188").Value = 100112028
$ws.Range("G188").Value = "Sandia"
$ws.Range("H188").Value = "Sin especificar"
$ws.Range("I188").Value = "Tercera"
$ws.Range("J188").Value = 3000
$ws.Range("K188").Value = 2000
$ws.Range("L188").Value = 2000
$ws.Range("M188").Value = 2000
$ws.Range("N188").Value = "$/unidad"
$ws.Range("O188").Value = "Región del Maule"
$ws.Range("P188").Value = 2000
$ws.Range("Q188").Value = 1
$ws.Range("R188").Value = "Hortaliza"

Write-Host "Done. Used range:" $ws.UsedRange.Address()
